$wb = $excel.ActiveWorkbook

# --- Remove the "DeleteList" sheet and rename "RegisterList" to "RequestList" ---
$wb.Worksheets("DeleteList").Delete() | Out-Null
$wb.Worksheets("RegisterList").Name = "RequestList"

# --- Insert a new header row at the top of "UserList" showing the column parameters ---
$ws = $wb.Worksheets("UserList")
$ws.Rows("1:1").Insert() | Out-Null

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "DoB"
$ws.Range("E1").Value = "Card #"
$ws.Range("F1").Value = "Last Accessed"
$ws.Range("G1").Value = "Employee Status"
$ws.Range("H1").Value = "Password"

# Match the selection left behind after editing
$ws.Range("B5").Select() | Out-Null
